# Apply the 2025-08-24 daily data update to the CTA violent crime YTD workbook.
# For each affected worksheet, update the listed cells to their new values.

$wb = $excel.ActiveWorkbook

$changes = @{
  "Citywide Totals" = @{
    "H2" = 70
    "D3" = 93
    "L3" = 159
    "B6" = 251
    "C6" = 321
    "E6" = 284
    "F6" = 375
    "G6" = 327
    "H6" = 291
    "I6" = 362
    "J6" = 272
    "K6" = 345
    "L6" = 309
    "B7" = 343
    "C7" = 432
    "D7" = 450
    "E7" = 433
    "F7" = 532
    "G7" = 477
    "H7" = 462
    "I7" = 594
    "J7" = 507
    "K7" = 609
    "L7" = 580
  }
  "Garfield Park" = @{
    "B6" = 9
    "I6" = 23
    "B7" = 11
    "I7" = 35
  }
  "Grand Crossing" = @{
    "C6" = 23
    "K6" = 31
    "C7" = 28
    "K7" = 47
  }
  "Humboldt Park" = @{
    "J5" = 6
    "J6" = 10
  }
  "By Neighborhood" = @{
    "F8" = 34
    "D19" = 19
    "B32" = 11
    "I32" = 35
    "C36" = 28
    "K36" = 47
    "J41" = 10
    "B53" = 33
    "C53" = 38
    "E53" = 52
    "G53" = 60
    "H53" = 55
    "I53" = 93
    "L54" = 2
    "G76" = 10
    "J76" = 11
    "L76" = 15
    "I77" = 32
    "J77" = 23
    "C85" = 12
    "L94" = 11
    "B98" = 343
    "C98" = 432
    "D98" = 450
    "E98" = 433
    "F98" = 532
    "G98" = 477
    "H98" = 462
    "I98" = 594
    "J98" = 507
    "K98" = 609
    "L98" = 580
  }
  "Loop" = @{
    "H2" = 6
    "B6" = 20
    "C6" = 23
    "E6" = 41
    "G6" = 41
    "H6" = 38
    "I6" = 60
    "B7" = 33
    "C7" = 38
    "E7" = 52
    "G7" = 60
    "H7" = 55
    "I7" = 93
  }
  "Rogers Park" = @{
    "G5" = 5
    "J5" = 9
    "L5" = 9
    "G6" = 10
    "J6" = 11
    "L6" = 15
  }
  "Roseland" = @{
    "I6" = 20
    "J6" = 13
    "I7" = 32
    "J7" = 23
  }
  "Chatham" = @{
    "D3" = 4
    "D6" = 19
  }
  "United Center" = @{
    "C4" = 9
    "C5" = 12
  }
  "Lower West Side" = @{
    "L3" = 1
    "L5" = 2
  }
  "West Town" = @{
    "L5" = 4
    "L6" = 11
  }
  "Austin" = @{
    "F5" = 23
    "F6" = 34
  }
}

foreach ($sheetName in $changes.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellMap = $changes[$sheetName]
    foreach ($cellRef in $cellMap.Keys) {
        $ws.Range($cellRef).Value = $cellMap[$cellRef]
    }
}
